$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Transmitance values (column B, rows 3-18) from 1 to 100
# to account for the extra ordinary ray for polarimetric acquisitions
$ws.Range("B3:B18").Value = 100

# Update the active cell / selection to B18
$ws.Range("B18").Select()
